$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Slit2", "Robo2", "ECs", 3, 1, 0.176022, 0.5280659999999999, 0.03293066697281707, 0.03293066697281707, 3, 1, 0.9052683333333333, 2.715805, 0.5200140314301739, 0.5200140314301738, 0.15934714257, 1.43412428313, 0.01712440889021909, 0.01712440889021908),
    @("ECs", "Slit2", "Robo2", "FAPs", 3, 1, 0.176022, 0.5280659999999999, 0.03293066697281707, 0.03293066697281707, 3, 1, 0.7756663333333332, 2.326999, 0.4455666482402024, 0.4455666482402024, 0.136534339326, 1.228809053934, 0.01467280690739244, 0.01467280690739243),
    @("ECs", "Slit2", "Robo2", "MuSCs", 3, 1, 0.176022, 0.5280659999999999, 0.03293066697281707, 0.03293066697281707, 2, 0.6666666666666666, 0.04422333333333334, 0.13267, 0.02540324564902162, 0.02540324564902162, 0.00778427958, 0.07005851621999999, 0.0008365458224965954, 0.0008365458224965951),
    @("ECs", "Slit2", "Robo2", "Resolving-Mac", 3, 1, 0.176022, 0.5280659999999999, 0.03293066697281707, 0.03293066697281707, 1, 0.3333333333333333, 0.01569566666666667, 0.047087, 0.009016074680602103, 0.009016074680602101, 0.002762782638, 0.024865043742, 0.0002969053527089559, 0.0002969053527089558),
    @("FAPs", "Slit2", "Robo2", "ECs", 3, 1, 4.277274333333334, 12.831823, 0.8002039325901205, 0.8002039325901203, 3, 1, 0.9052683333333333, 2.715805, 0.5200140314301739, 0.5200140314301738, 3.872081006946111, 34.848729062515, 0.4161172729524676, 0.4161172729524675),
    @("FAPs", "Slit2", "Robo2", "FAPs", 3, 1, 4.277274333333334, 12.831823, 0.8002039325901205, 0.8002039325901203, 3, 1, 0.7756663333333332, 2.326999, 0.4455666482402024, 0.4455666482402024, 3.317737698797444, 29.859639289177, 0.3565441841528089, 0.3565441841528088),
    @("FAPs", "Slit2", "Robo2", "MuSCs", 3, 1, 4.277274333333334, 12.831823, 0.8002039325901205, 0.8002039325901203, 2, 0.6666666666666666, 0.04422333333333334, 0.13267, 0.02540324564902162, 0.02540324564902162, 0.1891553286011111, 1.70239795741, 0.02032777706889997, 0.02032777706889996),
    @("FAPs", "Slit2", "Robo2", "Resolving-Mac", 3, 1, 4.277274333333334, 12.831823, 0.8002039325901205, 0.8002039325901203, 1, 0.3333333333333333, 0.01569566666666667, 0.047087, 0.009016074680602103, 0.009016074680602101, 0.06713467217788889, 0.604212049601, 0.007214698415944017, 0.007214698415944015),
    @("MuSCs", "Slit2", "Robo2", "ECs", 3, 1, 0.891934, 2.675802, 0.1668654004370625, 0.1668654004370625, 3, 1, 0.9052683333333333, 2.715805, 0.5200140314301739, 0.5200140314301738, 0.8074396056233334, 7.26695645061, 0.0867723495874872, 0.08677234958748718),
    @("MuSCs", "Slit2", "Robo2", "FAPs", 3, 1, 0.891934, 2.675802, 0.1668654004370625, 0.1668654004370625, 3, 1, 0.7756663333333332, 2.326999, 0.4455666482402024, 0.4455666482402024, 0.6918431753553332, 6.226588578197999, 0.07434965718000117, 0.07434965718000117),
    @("MuSCs", "Slit2", "Robo2", "MuSCs", 3, 1, 0.891934, 2.675802, 0.1668654004370625, 0.1668654004370625, 2, 0.6666666666666666, 0.04422333333333334, 0.13267, 0.02540324564902162, 0.02540324564902162, 0.03944429459333334, 0.35499865134, 0.00423892275762506, 0.004238922757625059),
    @("MuSCs", "Slit2", "Robo2", "Resolving-Mac", 3, 1, 0.891934, 2.675802, 0.1668654004370625, 0.1668654004370625, 1, 0.3333333333333333, 0.01569566666666667, 0.047087, 0.009016074680602103, 0.009016074680602101, 0.01399949875266667, 0.125995488774, 0.001504470911949131, 0.001504470911949131)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowData[$j]
    }
}
